$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (2025-09-... / serial 45556): mark habits done, add "play basketball" as the
# new habit for the week in column C
$ws.Range("B6").Value = "done"
$ws.Range("C6").Value = "play basketball"
$ws.Range("D6").Value = "done"

# Row 7 (serial 45557): same pattern
$ws.Range("B7").Value = "done"
$ws.Range("C7").Value = "play basketball"
$ws.Range("D7").Value = "done"

# Move the active selection to D7, matching where editing left off
$ws.Range("D7").Select()
